$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting
# (so values like "1.00" or "8.00" do not get auto-converted to numbers).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D13", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D35", "D39", "D40", "D41", "D42", "D45", "D47", "D49")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.924.91"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "3.371.66"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "569.91"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "136.29"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.370.60"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("D13").Value = "3.946.70"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "25.61"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "3.375.72"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "61.095.62"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "9.32"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "373.84"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "3.516.05"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "70.98"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  +11.19%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.63"
$ws.Range("E29").Value = "  -6.13%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "23.27"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "164.75"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "0.0755"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "0.773"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").Value = "4.31"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "2.533.26"
$ws.Range("E47").Value = "  +8.44%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "22.78"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("E51").Value = "  -1.44%  "
